$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 22

$ws.Range("C2").Select()
